$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.725.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "'1.852.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").Value = "'312.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("D6").Value = "'1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("D7").Value = "'0.4287"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.00%  "
$ws.Range("D8").Value = "'0.3593"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.40%  "
$ws.Range("D9").Value = "'0.07313"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("D10").Value = "'0.8739"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.77%  "
$ws.Range("D11").Value = "'20.79"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("D12").Value = "'1.926.02"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.41%  "
$ws.Range("D13").Value = "'6.558"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("D14").Value = "'5.340"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("D15").Value = "'0.07013"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.93%  "
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").Value = "'79.72"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").Value = "'0.000008957"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("D19").Value = "'1.003"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").Value = "'15.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.72%  "
$ws.Range("D21").Value = "'27.634.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").Value = "'10.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.84%  "
$ws.Range("D24").Value = "'2.053.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").Value = "'1.995"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.46%  "
$ws.Range("D26").Value = "'155.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.90%  "
$ws.Range("D27").Value = "'18.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.35%  "
$ws.Range("D28").Value = "'120.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.61%  "
$ws.Range("D29").Value = "'5.282"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").Value = "'1.882"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.10%  "
$ws.Range("D31").Value = "'0.08925"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("D32").Value = "'0.7581"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.80%  "
$ws.Range("D33").Value = "'2.975"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.18%  "
$ws.Range("D34").Value = "'4.518"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.36%  "
$ws.Range("D35").Value = "'1.129"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.92%  "
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("D37").Value = "'0.05432"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.08%  "
$ws.Range("D38").Value = "'1.101"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("D40").Value = "'2.833"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("D41").Value = "'0.1669"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("D42").Value = "'0.5070"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("D43").Value = "'6.633"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.82%  "
$ws.Range("D44").Value = "'8.430"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.74%  "
$ws.Range("D45").Value = "'106.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.46%  "
$ws.Range("D46").Value = "'0.06532"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.00%  "
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("D48").Value = "'0.4676"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("D49").Value = "'1.005"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.53%  "
$ws.Range("E50").Value = "  -0.52%  "
$ws.Range("D51").Value = "'1.796"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.05%  "
